# Adds a new block of 4 rows (model_s_rep = 20, i.e. "S = 20") to the DML
# binary-estimation results sheet, mirroring the existing per-S blocks
# (APO_0 / ATE / ATTE / APO_1) already present for S = 5, 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54 ---
$ws.Cells.Item(54, 1).Value = 'grade'
$ws.Cells.Item(54, 2).Value = 'controls_same_outcome'
$ws.Cells.Item(54, 3).Value = 'down'
$ws.Cells.Item(54, 4).Value = 'weekly'
$ws.Cells.Item(54, 5).Value = 'yes'
$ws.Cells.Item(54, 6).Value = 'all'
$ws.Cells.Item(54, 7).Value = 'lasso'
$ws.Cells.Item(54, 8).Value = 4.0
$ws.Cells.Item(54, 9).Value = 4.0
$ws.Cells.Item(54, 10).Value = 20.0
$ws.Cells.Item(54, 11).NumberFormat = "@"
$ws.Cells.Item(54, 11).Value = '0.01'
$ws.Cells.Item(54, 11).ClearFormats()
$ws.Cells.Item(54, 12).Value = 'no_treatment_outcome_lags'
$ws.Cells.Item(54, 13).Value = 'yes'
$ws.Cells.Item(54, 14).Value = 11122.0
$ws.Cells.Item(54, 15).Value = 11117.0
$ws.Cells.Item(54, 16).Value = 29.3275
$ws.Cells.Item(54, 17).Value = 81.3
$ws.Cells.Item(54, 18).Value = 50.46
$ws.Cells.Item(54, 19).Value = 'no'
$ws.Cells.Item(54, 20).Value = 'APO_0'
$ws.Cells.Item(54, 21).Value = 0.006448414444749802
$ws.Cells.Item(54, 22).Value = 0.0063458116204920915
$ws.Cells.Item(54, 23).Value = 0.00586027410042029
$ws.Cells.Item(54, 24).Value = 0.004210494688011051
$ws.Cells.Item(54, 25).Value = 1.1432440603525185
$ws.Cells.Item(54, 26).Value = 1.5285777332213555
$ws.Cells.Item(54, 27).Value = 0.27229161339660796
$ws.Cells.Item(54, 28).Value = 0.1263975975859819
$ws.Cells.Item(54, 29).Value = 0.006430679805759445
$ws.Cells.Item(54, 30).Value = 0.006466149083740158
$ws.Cells.Item(54, 31).Value = 0.006333069622083291
$ws.Cells.Item(54, 32).Value = 0.006358553618900892
$ws.Cells.Item(54, 33).Value = 0.8187340963852295
$ws.Cells.Item(54, 34).Value = 0.18039209516533683
$ws.Cells.Item(54, 35).Value = 0.5636796091397615
$ws.Cells.Item(54, 36).Value = 0.7394059918206516
$ws.Cells.Item(54, 37).Value = 0.4907011215366788
$ws.Cells.Item(54, 38).Value = 0.5181675881025801
$ws.Cells.Item(54, 39).Value = 251.09413390640668
$ws.Cells.Item(54, 40).Value = 269.898349306286
$ws.Cells.Item(54, 41).Value = 0.41150056245865535
$ws.Cells.Item(54, 42).Value = 0.45090114000864906
$ws.Cells.Item(54, 43).Value = 0.6413301610949432
$ws.Cells.Item(54, 44).Value = 0.670904004914734
$ws.Cells.Item(54, 45).Value = '16.4780809044176 hours'
$ws.Cells.Item(54, 46).Value = '2023-04-04 00:46:31'

# --- Row 55 ---
$ws.Cells.Item(55, 1).Value = 'grade'
$ws.Cells.Item(55, 2).Value = 'controls_same_outcome'
$ws.Cells.Item(55, 3).Value = 'down'
$ws.Cells.Item(55, 4).Value = 'weekly'
$ws.Cells.Item(55, 5).Value = 'yes'
$ws.Cells.Item(55, 6).Value = 'all'
$ws.Cells.Item(55, 7).Value = 'lasso'
$ws.Cells.Item(55, 8).Value = 4.0
$ws.Cells.Item(55, 9).Value = 4.0
$ws.Cells.Item(55, 10).Value = 20.0
$ws.Cells.Item(55, 11).NumberFormat = "@"
$ws.Cells.Item(55, 11).Value = '0.01'
$ws.Cells.Item(55, 11).ClearFormats()
$ws.Cells.Item(55, 12).Value = 'no_treatment_outcome_lags'
$ws.Cells.Item(55, 13).Value = 'yes'
$ws.Cells.Item(55, 14).Value = 11122.0
$ws.Cells.Item(55, 15).Value = 11117.0
$ws.Cells.Item(55, 16).Value = 29.3275
$ws.Cells.Item(55, 17).Value = 81.3
$ws.Cells.Item(55, 18).Value = 50.46
$ws.Cells.Item(55, 19).Value = 'no_yes'
$ws.Cells.Item(55, 20).Value = 'ATE'
$ws.Cells.Item(55, 21).Value = -0.03250071298585326
$ws.Cells.Item(55, 22).Value = -0.033667872232583475
$ws.Cells.Item(55, 23).Value = 0.016275683729889523
$ws.Cells.Item(55, 24).Value = 0.016208876818408095
$ws.Cells.Item(55, 25).Value = -1.9972308491074615
$ws.Cells.Item(55, 26).Value = -2.061890683621395
$ws.Cells.Item(55, 27).Value = 0.04714637619472564
$ws.Cells.Item(55, 28).Value = 0.039241323127657796
$ws.Cells.Item(55, 29).Value = -0.03254996723067606
$ws.Cells.Item(55, 30).Value = -0.03245145874103046
$ws.Cells.Item(55, 31).Value = -0.033716924303167174
$ws.Cells.Item(55, 32).Value = -0.033618820161999775
$ws.Cells.Item(55, 33).Value = 0.8187340963852295
$ws.Cells.Item(55, 34).Value = 0.18039209516533683
$ws.Cells.Item(55, 35).Value = 0.5636796091397615
$ws.Cells.Item(55, 36).Value = 0.7394059918206516
$ws.Cells.Item(55, 37).Value = 0.4907011215366788
$ws.Cells.Item(55, 38).Value = 0.5181675881025801
$ws.Cells.Item(55, 39).Value = 251.09413390640668
$ws.Cells.Item(55, 40).Value = 269.898349306286
$ws.Cells.Item(55, 41).Value = 0.41150056245865535
$ws.Cells.Item(55, 42).Value = 0.45090114000864906
$ws.Cells.Item(55, 43).Value = 0.6413301610949432
$ws.Cells.Item(55, 44).Value = 0.670904004914734
$ws.Cells.Item(55, 45).Value = '16.4780811536312 hours'
$ws.Cells.Item(55, 46).Value = '2023-04-04 00:46:31'

# --- Row 56 ---
$ws.Cells.Item(56, 1).Value = 'grade'
$ws.Cells.Item(56, 2).Value = 'controls_same_outcome'
$ws.Cells.Item(56, 3).Value = 'down'
$ws.Cells.Item(56, 4).Value = 'weekly'
$ws.Cells.Item(56, 5).Value = 'yes'
$ws.Cells.Item(56, 6).Value = 'all'
$ws.Cells.Item(56, 7).Value = 'lasso'
$ws.Cells.Item(56, 8).Value = 4.0
$ws.Cells.Item(56, 9).Value = 4.0
$ws.Cells.Item(56, 10).Value = 20.0
$ws.Cells.Item(56, 11).NumberFormat = "@"
$ws.Cells.Item(56, 11).Value = '0.01'
$ws.Cells.Item(56, 11).ClearFormats()
$ws.Cells.Item(56, 12).Value = 'no_treatment_outcome_lags'
$ws.Cells.Item(56, 13).Value = 'yes'
$ws.Cells.Item(56, 14).Value = 11122.0
$ws.Cells.Item(56, 15).Value = 11117.0
$ws.Cells.Item(56, 16).Value = 29.3275
$ws.Cells.Item(56, 17).Value = 81.3
$ws.Cells.Item(56, 18).Value = 50.46
$ws.Cells.Item(56, 19).Value = 'no_yes'
$ws.Cells.Item(56, 20).Value = 'ATTE'
$ws.Cells.Item(56, 21).Value = -0.03172562267050359
$ws.Cells.Item(56, 22).Value = -0.031737564722237696
$ws.Cells.Item(56, 23).Value = 0.01379196527947356
$ws.Cells.Item(56, 24).Value = 0.013800067207272749
$ws.Cells.Item(56, 25).Value = -2.3003119953999787
$ws.Cells.Item(56, 26).Value = -2.305891994179437
$ws.Cells.Item(56, 27).Value = 0.022750711970795257
$ws.Cells.Item(56, 28).Value = 0.02113497066269257
$ws.Cells.Item(56, 29).Value = -0.031767360568937635
$ws.Cells.Item(56, 30).Value = -0.03168388477206955
$ws.Cells.Item(56, 31).Value = -0.031779327139109304
$ws.Cells.Item(56, 32).Value = -0.03169580230536609
$ws.Cells.Item(56, 33).Value = 0.8187340963852295
$ws.Cells.Item(56, 34).Value = 0.18039209516533683
$ws.Cells.Item(56, 35).Value = 0.5636796091397615
$ws.Cells.Item(56, 36).Value = 0.7394059918206516
$ws.Cells.Item(56, 37).Value = 0.4907011215366788
$ws.Cells.Item(56, 38).Value = 0.5181675881025801
$ws.Cells.Item(56, 39).Value = 251.09413390640668
$ws.Cells.Item(56, 40).Value = 269.898349306286
$ws.Cells.Item(56, 41).Value = 0.41150056245865535
$ws.Cells.Item(56, 42).Value = 0.45090114000864906
$ws.Cells.Item(56, 43).Value = 0.6413301610949432
$ws.Cells.Item(56, 44).Value = 0.670904004914734
$ws.Cells.Item(56, 45).Value = '16.4780811536312 hours'
$ws.Cells.Item(56, 46).Value = '2023-04-04 00:46:31'

# --- Row 57 ---
$ws.Cells.Item(57, 1).Value = 'grade'
$ws.Cells.Item(57, 2).Value = 'controls_same_outcome'
$ws.Cells.Item(57, 3).Value = 'down'
$ws.Cells.Item(57, 4).Value = 'weekly'
$ws.Cells.Item(57, 5).Value = 'yes'
$ws.Cells.Item(57, 6).Value = 'all'
$ws.Cells.Item(57, 7).Value = 'lasso'
$ws.Cells.Item(57, 8).Value = 4.0
$ws.Cells.Item(57, 9).Value = 4.0
$ws.Cells.Item(57, 10).Value = 20.0
$ws.Cells.Item(57, 11).NumberFormat = "@"
$ws.Cells.Item(57, 11).Value = '0.01'
$ws.Cells.Item(57, 11).ClearFormats()
$ws.Cells.Item(57, 12).Value = 'no_treatment_outcome_lags'
$ws.Cells.Item(57, 13).Value = 'yes'
$ws.Cells.Item(57, 14).Value = 11122.0
$ws.Cells.Item(57, 15).Value = 11117.0
$ws.Cells.Item(57, 16).Value = 29.3275
$ws.Cells.Item(57, 17).Value = 81.3
$ws.Cells.Item(57, 18).Value = 50.46
$ws.Cells.Item(57, 19).Value = 'yes'
$ws.Cells.Item(57, 20).Value = 'APO_1'
$ws.Cells.Item(57, 21).Value = -0.02605229854110346
$ws.Cells.Item(57, 22).Value = -0.0265058312944634
$ws.Cells.Item(57, 23).Value = 0.007956742527413576
$ws.Cells.Item(57, 24).Value = 0.005445171150942421
$ws.Cells.Item(57, 25).Value = -3.3500071697880087
$ws.Cells.Item(57, 26).Value = -4.5560050377228425
$ws.Cells.Item(57, 27).Value = 0.0019845663800074297
$ws.Cells.Item(57, 28).Value = [double]"5.269068479047221E-6"
$ws.Cells.Item(57, 29).Value = -0.026076377612123836
$ws.Cells.Item(57, 30).Value = -0.026028219470083082
$ws.Cells.Item(57, 31).Value = -0.026522309729251244
$ws.Cells.Item(57, 32).Value = -0.026489352859675554
$ws.Cells.Item(57, 33).Value = 0.8187340963852295
$ws.Cells.Item(57, 34).Value = 0.18039209516533683
$ws.Cells.Item(57, 35).Value = 0.5636796091397615
$ws.Cells.Item(57, 36).Value = 0.7394059918206516
$ws.Cells.Item(57, 37).Value = 0.4907011215366788
$ws.Cells.Item(57, 38).Value = 0.5181675881025801
$ws.Cells.Item(57, 39).Value = 251.09413390640668
$ws.Cells.Item(57, 40).Value = 269.898349306286
$ws.Cells.Item(57, 41).Value = 0.41150056245865535
$ws.Cells.Item(57, 42).Value = 0.45090114000864906
$ws.Cells.Item(57, 43).Value = 0.6413301610949432
$ws.Cells.Item(57, 44).Value = 0.670904004914734
$ws.Cells.Item(57, 45).Value = '16.4780811536312 hours'
$ws.Cells.Item(57, 46).Value = '2023-04-04 00:46:31'

